$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up formatting on the existing row 3 email cells so they match the
#     same look already used on row 2 (C2/D2), instead of the generic
#     "hyperlink placeholder" style they had before.
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C3").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Add the new review row (row 4): helix app review from rocketaso@gmail.com
$ws.Range("A4").Value = "com.singleton.helix"
$ws.Range("B4").Value = "helix"
$ws.Range("C4").Value = "rocketaso@gmail.com"
$ws.Range("E4").Value = "27/5/2019 15:56"
$ws.Range("F4").Value = "great music and fun. Love to spend hours in it."

# Match formatting of row 4 to its neighbouring rows.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- New hyperlink for the added email cell.
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:rocketaso@gmail.com", [Type]::Missing, [Type]::Missing, "rocketaso@gmail.com") | Out-Null

# --- Update the active selection.
$ws.Range("F4").Select() | Out-Null

Write-Host "edit complete"
